$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the course rows: row2 -> LBRL750 data, row3 -> LBRL751 data, row4 -> LBRL752 data
$ws.Range("A2").Value = "LBRL750"
$ws.Range("D2").Value = "LBRL 750 - Co-op I"

$ws.Range("A3").Value = "LBRL751"
$ws.Range("D3").Value = "LBRL 751 - Co-op II"

$ws.Range("A4").Value = "LBRL752"
$ws.Range("D4").Value = "LBRL 752 - Co-op III"

# Update the active cell / selection
$ws.Range("G7").Select()
